$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values stay plain text (several look like decimal numbers and
# would otherwise be auto-converted to numeric values by Excel).

# Row 38 <-> Row 39 swap (ImmutableX now ranked above VeChain) plus updated values
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.574'
$ws.Range("E38").Value = '  -5.29%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0172'
$ws.Range("E39").Value = '  -1.70%  '

# Updated Price (D) and Volume(1h) (E) figures for the remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.596.20'
$ws.Range("E2").Value = '  -1.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.666.41'
$ws.Range("E3").Value = '  -3.25%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.97'
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.508'
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.81'
$ws.Range("E8").Value = '  -1.96%  '
$ws.Range("E9").Value = '  -0.75%  '
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.902.68'
$ws.Range("E12").Value = '  -3.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.655.52'
$ws.Range("E13").Value = '  -3.87%  '
$ws.Range("E14").Value = '  -3.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.558'
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.29'
$ws.Range("E16").Value = '  -1.83%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.579.31'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '243.06'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0729'
$ws.Range("E19").Value = '  -3.33%  '
$ws.Range("E20").Value = '  -4.63%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("E22").Value = '  -3.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.30'
$ws.Range("E23").Value = '  -3.78%  '
$ws.Range("E24").Value = '  -4.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.99'
$ws.Range("E25").Value = '  -1.12%  '
$ws.Range("E26").Value = '  -3.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.45'
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("E28").Value = '  +0.09%  '
$ws.Range("E29").Value = '  -2.43%  '
$ws.Range("E30").Value = '  +2.98%  '
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("E32").Value = '  -2.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.468.98'
$ws.Range("E33").Value = '  -1.48%  '
$ws.Range("E34").Value = '  -4.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  -5.15%  '
$ws.Range("E36").Value = '  -1.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.928'
$ws.Range("E37").Value = '  -2.75%  '
$ws.Range("E40").Value = '  -1.60%  '
$ws.Range("E41").Value = '  -4.84%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("E43").Value = '  -3.07%  '
$ws.Range("E44").Value = '  -7.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.809.59'
$ws.Range("E45").Value = '  -3.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.786'
$ws.Range("E46").Value = '  -1.69%  '
$ws.Range("E47").Value = '  -2.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.33'
$ws.Range("E48").Value = '  -1.81%  '
$ws.Range("E49").Value = '  -4.08%  '
$ws.Range("E50").Value = '  -2.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.87'
$ws.Range("E51").Value = '  -4.70%  '
